# LS_Trafo is added v1
# Mark the "Boad design for analog circuit and reduce isolation" row (A6) as
# Urgent, matching the styling already used on A7 (red font / highlighted fill).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Copy()
$ws.Range("A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Remove the stray test content that had been added to A8, A9 and A11,
# leaving the rows themselves (and their existing formatting) in place.
$ws.Range("A8").Value = ""
$ws.Range("A9").Value = ""
$ws.Range("A11").Value = ""

# Keep rows 8, 9 and 11 present (with their recalculated height) even though
# they no longer hold any text.
$ws.Rows.Item(8).RowHeight = 14.45
$ws.Rows.Item(9).RowHeight = 14.45
$ws.Rows.Item(11).RowHeight = 14.45
